# Apply "graphs rh-yield and fix stats" edits to the stats sheet.
# Rows 6-13 get re-ordered/corrected (the 2020/11-12 and 2021/11-12 blocks
# swap position, and several mean/max values are corrected), and rows 15-16
# get corrected mean_Y / max_Y values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 6 through 16: Year, Month, mean_Y, max_Y, min_Y
$data = @{
    6  = @(2021, 11, 0,                   0,            0)
    7  = @(2021, 12, 0.190827956989248,   0.471,        0)
    8  = @(2022, 1,  0.22780376344086,    0.454,        0)
    9  = @(2022, 2,  0.0550744047619051,  0.397,        0)
    10 = @(2020, 11, 0.0033129359125,     0.169230769,  0)
    11 = @(2020, 12, 0.229241229139785,   0.458,        0)
    12 = @(2021, 1,  0.238112903225806,   0.471,        0)
    13 = @(2021, 2,  0.0238303571428571,  0.207,        0)
    14 = @(2022, 11, 0.00167083333333335, 0.176,        0)
    15 = @(2022, 12, 0.157638440860215,   0.8,          0)
    16 = @(2023, 1,  0.119631720430108,   0.714,        0)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}
